$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.09569222968600938
$ws.Range("E2").Value = 0.001569214846296548
$ws.Range("D3").Value = 0.03650567731802929
$ws.Range("E3").Value = 0.0009654804221026359
$ws.Range("D4").Value = 0.07255159012043401
$ws.Range("E4").Value = 0.001282918626960443
$ws.Range("B5").Value = 2542.458447620835
$ws.Range("D5").Value = 0.01396910341455489
$ws.Range("E5").Value = 0.0008497692244262805
$ws.Range("B6").Value = 2551.877524835181
$ws.Range("D6").Value = 0.01391083285112843
$ws.Range("E6").Value = 0.0007044601341936426
$ws.Range("B7").Value = 2565.405862269178
$ws.Range("D7").Value = 0.01852630631962342
$ws.Range("E7").Value = 0.000736635272234543
$ws.Range("B8").Value = 2578.18433748
$ws.Range("D8").Value = 0.01754401233912021
$ws.Range("E8").Value = 0.0006722683748416727
$ws.Range("B9").Value = 2607.312573558408
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0.01131255283757036
$ws.Range("E9").Value = 0.0007706521218561054
$ws.Range("B10").Value = 2618.528684202817
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 0.0319262839589859
$ws.Range("E10").Value = 0.001257832424016735
$ws.Range("B11").Value = 2632.326939830645
$ws.Range("D11").Value = 0.003154718946369824
$ws.Range("E11").Value = 0.0006379809580535667
$ws.Range("B12").Value = 2641.177916256609
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 0.01169254794461194
$ws.Range("E12").Value = 0.0006787850309811989
$ws.Range("B13").Value = 2656.047601072758
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0.01497793041804329
$ws.Range("E13").Value = 0.0007699777771719565
$ws.Range("B14").Value = 2671.279693108128
$ws.Range("D14").Value = 0.008702213210834898
$ws.Range("E14").Value = 0.0007877476392629149
$ws.Range("B15").Value = 2702.631537854082
$ws.Range("D15").Value = 0.01828680639890687
$ws.Range("E15").Value = 0.0009455490873658192
$ws.Range("B16").Value = 2712.689927009923
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 0.03920612912783229
$ws.Range("E16").Value = 0.002095339442698875
$ws.Range("B17").Value = 2727.142436180641
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 0.01370147051829513
$ws.Range("E17").Value = 0.001016906015029717
$ws.Range("B18").Value = 2747.912861313825
$ws.Range("D18").Value = 0.03299248711461778
$ws.Range("E18").Value = 0.001309875133076821
$ws.Range("B19").Value = 2757.600106909915
$ws.Range("D19").Value = 0.01423550172483481
$ws.Range("E19").Value = 0.001036199070880823
$ws.Range("B20").Value = 2776.47647050831
$ws.Range("D20").Value = 0.01942991730488781
$ws.Range("E20").Value = 0.0007422957683502713
$ws.Range("D21").Value = 0.04939162516526741
$ws.Range("E21").Value = 0.001230611626192194
$ws.Range("B22").Value = 2826.876574034807
$ws.Range("D22").Value = 0.003268186357642634
$ws.Range("E22").Value = 0.0003100255833336887
$ws.Range("B23").Value = 2841.94965238447
$ws.Range("D23").Value = 0.001172823520964464
$ws.Range("E23").Value = 0.0002577634112009812
$ws.Range("B24").Value = 2869.393415639229
$ws.Range("D24").Value = 0.01767670067440815
$ws.Range("E24").Value = 0.0006782094441907681
$ws.Range("B25").Value = 2887.379387138704
$ws.Range("D25").Value = 0.01846547837100152
$ws.Range("E25").Value = 0.001648835835377327
$ws.Range("B26").Value = 2894.474014632406
$ws.Range("D26").Value = 0.01894322678320538
$ws.Range("E26").Value = 0.002991035807874604
$ws.Range("B27").Value = 2907.162093876396
$ws.Range("D27").Value = 0.01188704610445319
$ws.Range("E27").Value = 0.0008257352711212121
$ws.Range("B28").Value = 2920.055011583651
$ws.Range("D28").Value = 0.02474372606616985
$ws.Range("E28").Value = 0.002717519047729712
$ws.Range("B29").Value = 2933.320381360995
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0.02914793126332345
$ws.Range("E29").Value = 0.002871717365844667
$ws.Range("B30").Value = 2957.044063330045
$ws.Range("D30").Value = 0.03085275707638647
$ws.Range("E30").Value = 0.0008236471310576636
$ws.Range("B31").Value = 2980.865171207523
$ws.Range("C31").Value = 4
$ws.Range("D31").Value = 0.007898475902520222
$ws.Range("E31").Value = 0.0008585768435644596
$ws.Range("B32").Value = 2996.615218302426
$ws.Range("D32").Value = 0.02085218887280283
$ws.Range("E32").Value = 0.00422876557560344
$ws.Range("B33").Value = 3005.069066343128
$ws.Range("D33").Value = 0.01772861638531403
$ws.Range("E33").Value = 0.001234512602392027
$ws.Range("D34").Value = 0.0215692567624347
$ws.Range("E34").Value = 0.0008971119157286981
$ws.Range("B35").Value = 3082.097723639147
$ws.Range("D35").Value = 0.01862608143692363
$ws.Range("E35").Value = 0.003427198984394019
$ws.Range("B36").Value = 3096.949709226105
$ws.Range("D36").Value = 0.01058526044356712
$ws.Range("E36").Value = 0.0005906108348060491
$ws.Range("B37").Value = 3134.848126147962
$ws.Range("C37").Value = 3
$ws.Range("D37").Value = 0.006352551883367627
$ws.Range("E37").Value = 0.0004693479703699076
$ws.Range("B38").Value = 3158.424062319133
$ws.Range("D38").Value = 0.008156029184079714
$ws.Range("E38").Value = 0.002265563662244454
$ws.Range("B39").Value = 3170.520162281313
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = 0.003390170246613524
$ws.Range("E39").Value = 0.0005752274425661935
$ws.Range("B40").Value = 3182.783361911292
$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 0.005863079944031258
$ws.Range("E40").Value = 0.0004641604955691411
$ws.Range("B41").Value = 3189.946648598499
$ws.Range("D41").Value = 0.01243133923676314
$ws.Range("E41").Value = 0.0008052707941379539
$ws.Range("B42").Value = 3225.325087674019
$ws.Range("D42").Value = 0.00279989264324885
$ws.Range("E42").Value = 0.0005090713896816091
$ws.Range("B43").Value = 3245.046605571625
$ws.Range("D43").Value = 0.002610590011497431
$ws.Range("E43").Value = 0.0003625819460413099
$ws.Range("B44").Value = 3261.235108125469
$ws.Range("D44").Value = 0.006161086105706273
$ws.Range("E44").Value = 0.0005303384975836233
$ws.Range("B45").Value = 3278.01174742576
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 0.006806949329930652
$ws.Range("E45").Value = 0.000486301393733876
$ws.Range("B46").Value = 3299.742583620296
$ws.Range("D46").Value = 0.01142098472262458
$ws.Range("E46").Value = 0.0004708111657760201
$ws.Range("B47").Value = 3322.424071527266
$ws.Range("D47").Value = 0.004091902478616681
$ws.Range("E47").Value = 0.0003147866702038125
$ws.Range("B48").Value = 3343.020580858651
$ws.Range("D48").Value = 0.00882307035001133
$ws.Range("E48").Value = 0.0003878720360807331
$ws.Range("B49").Value = 3365.876821333213
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = 0.01115817986354503
$ws.Range("E49").Value = 0.000471422230556489
$ws.Range("B50").Value = 3384.081593992146
$ws.Range("D50").Value = 0.01344279858559126
$ws.Range("E50").Value = 0.0005672463968599721
$ws.Range("B51").Value = 3398.742210501987
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 0.006880520030706584
$ws.Range("E51").Value = 0.0005231769943595078
$ws.Range("B52").Value = 3415.445399821636
$ws.Range("D52").Value = 0.01344403256502429
$ws.Range("E52").Value = 0.0004760436909056625
$ws.Range("B53").Value = 3435.532051038655
$ws.Range("D53").Value = 0.008826142606088962
$ws.Range("E53").Value = 0.0003869125191954116
$ws.Range("B54").Value = 3457.480910751303
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = 0.01180867359787212
$ws.Range("E54").Value = 0.0004361614537595017
$ws.Range("B55").Value = 3476.652358107846
$ws.Range("D55").Value = 0.01001848923154773
$ws.Range("E55").Value = 0.002886683337903653
$ws.Range("B56").Value = 3491.023967613369
$ws.Range("D56").Value = 0.01054891832307162
$ws.Range("E56").Value = 0.0005235725708735842
$ws.Range("B57").Value = 3509.332337392685
$ws.Range("D57").Value = 0.02081597009038815
$ws.Range("E57").Value = 0.0007180112003547019
$ws.Range("B58").Value = 3523.826377073166
$ws.Range("D58").Value = 0.005217428418431261
$ws.Range("E58").Value = 0.0006623220357618579
$ws.Range("B59").Value = 3539.054419347468
$ws.Range("C59").Value = 3
$ws.Range("D59").Value = 0.003938370300913865
$ws.Range("E59").Value = 0.0006467206706962789
$ws.Range("B60").Value = 3552.499083567654
$ws.Range("C60").Value = 3
$ws.Range("D60").Value = 0.01984108787569523
$ws.Range("E60").Value = 0.0008044082145588651
$ws.Range("B61").Value = 3568.751028993931
$ws.Range("D61").Value = 0.009266643385596256
$ws.Range("E61").Value = 0.0006797471136224874
$ws.Range("B62").Value = 3579.273979331836
$ws.Range("D62").Value = 0.008432062153133716
$ws.Range("E62").Value = 0.0006958188838078449
$ws.Range("B63").Value = 3597.153147696692
$ws.Range("D63").Value = 0.04901155504291254
$ws.Range("E63").Value = 0.004633819749511721
